$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.00"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.88"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.383"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.390"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8154"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9343"

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "One"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.01121"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "8OneONEBestin24h"

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1432"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03520"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07404"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03055"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09410"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.007"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001600"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04810"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.005668"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.004160"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009917"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.667"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.424"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.00007002"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002901"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04015"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006401"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002901"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005931"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "43LocalTradersLCT"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005273"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9302"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002484"
